$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" '27.689.58'
Set-TextValue "E2" '  -0.71%  '
Set-TextValue "D3" '1.590.38'
Set-TextValue "E3" '  -2.40%  '
Set-TextValue "E4" '  +0.15%  '
Set-TextValue "D5" '208.52'
Set-TextValue "E5" '  -1.49%  '
Set-TextValue "E6" '  -3.11%  '
Set-TextValue "E7" '  +0.15%  '
Set-TextValue "D8" '22.25'
Set-TextValue "E8" '  -4.24%  '
Set-TextValue "E9" '  -1.86%  '
Set-TextValue "E10" '  -2.53%  '
Set-TextValue "D11" '0.0868'
Set-TextValue "E11" '  -1.44%  '
Set-TextValue "D12" '1.815.82'
Set-TextValue "E12" '  -2.37%  '
Set-TextValue "D13" '1.594.84'
Set-TextValue "E13" '  -2.31%  '
Set-TextValue "E15" '  -4.26%  '
Set-TextValue "D16" '63.51'
Set-TextValue "E16" '  -2.02%  '
Set-TextValue "D17" '27.658.34'
Set-TextValue "E17" '  -0.88%  '
Set-TextValue "D18" '219.91'
Set-TextValue "E18" '  -3.57%  '
Set-TextValue "D19" '0.0₃0696'
Set-TextValue "E19" '  -3.01%  '
Set-TextValue "E20" '  -3.62%  '
Set-TextValue "E21" '  +0.18%  '
Set-TextValue "E22" '  -4.61%  '
Set-TextValue "D23" '9.67'
Set-TextValue "E23" '  -2.96%  '
Set-TextValue "E24" '  -3.67%  '
Set-TextValue "D25" '153.91'
Set-TextValue "E25" '  -0.77%  '
Set-TextValue "D26" '6.82'
Set-TextValue "E26" '  -1.27%  '
Set-TextValue "E27" '  +0.17%  '
Set-TextValue "D28" '15.13'
Set-TextValue "E28" '  -1.86%  '
Set-TextValue "D29" '0.105'
Set-TextValue "E29" '  -4.74%  '
Set-TextValue "E30" '  -1.91%  '
Set-TextValue "D31" '0.0470'
Set-TextValue "E31" '  -2.21%  '
Set-TextValue "D32" '3.22'
Set-TextValue "E32" '  -5.21%  '
Set-TextValue "D33" '1.376.57'
Set-TextValue "E33" '  -2.58%  '
Set-TextValue "E34" '  -4.99%  '
Set-TextValue "D36" '0.973'
Set-TextValue "E36" '  -2.89%  '
Set-TextValue "E37" '  -0.48%  '
Set-TextValue "E38" '  -0.91%  '
Set-TextValue "E39" '  -2.93%  '
Set-TextValue "D40" '0.827'
Set-TextValue "E40" '  -2.80%  '
Set-TextValue "E41" '  +0.14%  '
Set-TextValue "E42" '  -3.73%  '
Set-TextValue "E43" '  -2.12%  '
Set-TextValue "E44" '  +2.30%  '
Set-TextValue "E45" '  -3.55%  '
Set-TextValue "D46" '1.73'
Set-TextValue "E46" '  -4.42%  '
Set-TextValue "D47" '1.726.32'
Set-TextValue "E47" '  -2.41%  '
Set-TextValue "D48" '87.24'
Set-TextValue "E48" '  -1.57%  '
Set-TextValue "E49" '  -0.86%  '
Set-TextValue "D50" '0.0965'
Set-TextValue "E50" '  -4.12%  '
Set-TextValue "E51" '  -1.46%  '
